$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("T2").Value = 2
$ws.Range("U2").Value = 2
$ws.Range("V2").Value = 0
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 0
$ws.Range("Y2").Value = 0
$ws.Range("Z2").Value = 0
$ws.Range("AA2").Value = 0
$ws.Range("AC2").Value = 3
$ws.Range("AD2").Value = "Диспаритет"

# Row 4
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 2
$ws.Range("W4").Value = 2
$ws.Range("X4").Value = 0
$ws.Range("Y4").Value = 0
$ws.Range("Z4").Value = 0
$ws.Range("AA4").Value = 0
$ws.Range("AC4").Value = 6

# Row 5
$ws.Range("T5").Value = 2
$ws.Range("U5").Value = 2
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 0
$ws.Range("X5").Value = 0
$ws.Range("Y5").Value = 0
$ws.Range("Z5").Value = 0
$ws.Range("AA5").Value = 0
$ws.Range("AC5").Value = 4
$ws.Range("AD5").Value = "Диспаритет"

# Row 6
$ws.Range("X6").Value = 0
$ws.Range("Y6").Value = 0
$ws.Range("AC6").Value = 0
$ws.Range("AD6").Value = "Диспаритет"

# Row 7
$ws.Range("V7").Value = 2
$ws.Range("W7").Value = 2
$ws.Range("AC7").Value = 9
